$wb = $excel.ActiveWorkbook

# --- Step 3 sheet: update reagent row (row 3) from "ammonium hydroxide solution" to "ammonia" ---
$ws3 = $wb.Worksheets.Item("Step 3")

# Update the molecule name and SMILES for the reagent
$ws3.Range("A3").Value = "ammonia"
$ws3.Range("B3").Value = "N"

# Update MW, buy CAD and buy mass values (G3/K3/L3 and the row2 totals recalc automatically)
$ws3.Range("C3").Value = 17.03
$ws3.Range("E3").Value = 828
$ws3.Range("F3").Value = 170

# Update the vendor hyperlink target (D3) to the new Sigma-Aldrich product page
$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://www.sigmaaldrich.com/catalog/product/aldrich/294993?lang=en&region=CA")
$ws3.Range("D3").Value = "https://www.sigmaaldrich.com/catalog/product/aldrich/294993?lang=en&region=CA"
$ws3.Range("D3").Style = "Hyperlink"

# --- Activate "Step 3" as the selected tab, with A3:F3 selected ---
$ws3.Activate()
$ws3.Range("A3:F3").Select()
